$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the workbook's window (workbookView minimized="1")
$wb.Windows.Item(1).WindowState = -4140

# Add the new note to cell A9 (creates new shared string entry)
$ws.Range("A9").Value = "speed and protection zone are in 570410 where 64 means that character is hasted and 16384 means that character is in pz 16384+64 means is in pz and hasted"

# Move the selection to A9 (as reflected in the diff's <selection activeCell="A9" sqref="A9"/>)
$ws.Range("A9").Select()
